$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / 1h-volume-change figures pulled in by
# the scheduled GitHub Actions scrape, including two rank ties that swapped
# order (LEO/LidoDAOToken and Stellar/WEMIXTOKEN).
#
# Price cells whose new text would otherwise look like a plain decimal
# number (e.g. "306.50", "0.9972") are forced to Text format first so
# Excel keeps them as the literal string instead of coercing them into a
# floating point value (which would also corrupt trailing zeros).

$ws.Range('D2').Value = '24.484.55'
$ws.Range('E2').Value = '  +9.28%  '
$ws.Range('D3').Value = '1.679.74'
$ws.Range('E3').Value = '  +4.94%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '306.50'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9972'
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3714'
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3443'
$ws.Range('E8').Value = '  +1.27%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.88'
$ws.Range('E9').Value = '  +11.76%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.182'
$ws.Range('E10').Value = '  +3.55%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07269'
$ws.Range('E11').Value = '  +3.15%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.9990'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.38'
$ws.Range('E13').Value = '  +3.57%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.125'
$ws.Range('E14').Value = '  +3.12%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.739'
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').Value = '1.677.25'
$ws.Range('E16').Value = '  +4.87%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001107'
$ws.Range('E17').Value = '  +2.13%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.9969'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06714'
$ws.Range('E19').Value = '  +0.82%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '81.27'
$ws.Range('E20').Value = '  +4.29%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.44'
$ws.Range('E21').Value = '  +1.99%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.097'
$ws.Range('E22').Value = '  +1.10%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.96'
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('D24').Value = '24.409.63'
$ws.Range('E24').Value = '  +8.74%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.427'
$ws.Range('E25').Value = '  +1.11%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.362'
$ws.Range('E26').Value = '  -12.04%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.662'
$ws.Range('E27').Value = '  +6.54%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '152.65'
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '19.55'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = '1.862.11'
$ws.Range('E30').Value = '  +4.82%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '126.69'
$ws.Range('E31').Value = '  +4.93%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.335'
$ws.Range('E32').Value = '  +4.72%  '
$ws.Range('E33').Value = '  -4.04%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.9680'
$ws.Range('E34').Value = '  +1.82%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.08465'
$ws.Range('E35').Value = '  +2.60%  '
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.689'
$ws.Range('E36').Value = '  +2.83%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '12.42'
$ws.Range('E37').Value = '  +4.73%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06488'
$ws.Range('E38').Value = '  +6.18%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.943'
$ws.Range('E39').Value = '  +3.92%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.342'
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.02339'
$ws.Range('E41').Value = '  +5.43%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.268'
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.2111'
$ws.Range('E43').Value = '  +4.09%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6176'
$ws.Range('E44').Value = '  +4.41%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.9970'
$ws.Range('E45').Value = '  +0.55%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.771'
$ws.Range('E46').Value = '  -2.12%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5942'
$ws.Range('E47').Value = '  +4.25%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '12.95'
$ws.Range('E48').Value = '  -1.51%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '126.96'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.030'
$ws.Range('E50').Value = '  +3.10%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07214'
$ws.Range('E51').Value = '  +5.85%  '

Write-Output "Applied 106 cell updates"
